# Update the vm_pu results table for Case_3_213 (res_bus) after
# switching the slack-bus voltage setpoint from 1.05 pu to 1.02 pu
# ("case with 380 kV done"). Rows 2-25 correspond to timesteps 0-23,
# columns B:F and I:N hold recomputed per-bus voltage magnitudes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030526003742771
$ws.Range("D2").Value = 1.040461522314663
$ws.Range("E2").Value = 1.030223834145874
$ws.Range("F2").Value = 1.051579515711766
$ws.Range("I2").Value = 1.040015485696387
$ws.Range("J2").Value = 1.035666763105735
$ws.Range("K2").Value = 1.043243804420683
$ws.Range("L2").Value = 1.033035466514235
$ws.Range("M2").Value = 1.054330625800739
$ws.Range("N2").Value = 1.015893170983726

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031292591898147
$ws.Range("D3").Value = 1.041066172200911
$ws.Range("E3").Value = 1.030870285676733
$ws.Range("F3").Value = 1.052355819993733
$ws.Range("I3").Value = 1.040209368943919
$ws.Range("J3").Value = 1.036075775591138
$ws.Range("K3").Value = 1.043659424005461
$ws.Range("L3").Value = 1.033490663163847
$ws.Range("M3").Value = 1.054919698084545
$ws.Range("N3").Value = 1.016029173046155

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031789277581261
$ws.Range("D4").Value = 1.041457958450094
$ws.Range("E4").Value = 1.031289534618939
$ws.Range("F4").Value = 1.052859063102972
$ws.Range("I4").Value = 1.040333837975304
$ws.Range("J4").Value = 1.036340398604276
$ws.Range("K4").Value = 1.04392819843977
$ws.Range("L4").Value = 1.033785470258638
$ws.Range("M4").Value = 1.055301127508581
$ws.Range("N4").Value = 1.016117143178089

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031998238399093
$ws.Range("D5").Value = 1.041622792084157
$ws.Range("E5").Value = 1.031466012708246
$ws.Range("F5").Value = 1.053070845121285
$ws.Range("I5").Value = 1.040385927707989
$ws.Range("J5").Value = 1.036451636310581
$ws.Range("K5").Value = 1.044041151668846
$ws.Range("L5").Value = 1.033909469003809
$ws.Range("M5").Value = 1.055461540927354
$ws.Range("N5").Value = 1.016154117617067

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032033332824297
$ws.Range("D6").Value = 1.041650475749431
$ws.Range("E6").Value = 1.0314956573415
$ws.Range("F6").Value = 1.053106417009967
$ws.Range("I6").Value = 1.040394659876452
$ws.Range("J6").Value = 1.036470312997076
$ws.Range("K6").Value = 1.044060114652788
$ws.Range("L6").Value = 1.033930292522694
$ws.Range("M6").Value = 1.055488478540388
$ws.Range("N6").Value = 1.016160325293198

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031792069122595
$ws.Range("D7").Value = 1.041460160469806
$ws.Range("E7").Value = 1.031291891842922
$ws.Range("F7").Value = 1.052861892088415
$ws.Range("I7").Value = 1.040334534933427
$ws.Range("J7").Value = 1.036341885008163
$ws.Range("K7").Value = 1.043929707883512
$ws.Range("L7").Value = 1.03378712689478
$ws.Range("M7").Value = 1.055303270725141
$ws.Range("N7").Value = 1.016117637265185

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03078493967525
$ws.Range("D8").Value = 1.04066575414587
$ws.Range("E8").Value = 1.030442107251185
$ws.Range("F8").Value = 1.051841679086738
$ws.Range("I8").Value = 1.040081213143788
$ws.Range("J8").Value = 1.035804997323672
$ws.Range("K8").Value = 1.043384297119188
$ws.Range("L8").Value = 1.033189246724411
$ws.Range("M8").Value = 1.054529650553775
$ws.Range("N8").Value = 1.015939139920757

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029015325954225
$ws.Range("D9").Value = 1.039270110066017
$ws.Range("E9").Value = 1.028952053425679
$ws.Range("F9").Value = 1.050051094303477
$ws.Range("I9").Value = 1.03962731316851
$ws.Range("J9").Value = 1.034858720197221
$ws.Range("K9").Value = 1.042422057096048
$ws.Range("L9").Value = 1.032137795234912
$ws.Range("M9").Value = 1.053168506745643
$ws.Range("N9").Value = 1.015624378325523

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027839110940941
$ws.Range("D10").Value = 1.038342621703502
$ws.Range("E10").Value = 1.027963760953374
$ws.Range("F10").Value = 1.048862316757851
$ws.Range("I10").Value = 1.039319714812407
$ws.Range("J10").Value = 1.034227805426091
$ws.Range("K10").Value = 1.041779871193852
$ws.Range("L10").Value = 1.031438319940078
$ws.Range("M10").Value = 1.052262580585823
$ws.Range("N10").Value = 1.015414413224361

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027330657622345
$ws.Range("D11").Value = 1.037941731223319
$ws.Range("E11").Value = 1.027537048592156
$ws.Range("F11").Value = 1.048348763777582
$ws.Range("I11").Value = 1.039185347200179
$ws.Range("J11").Value = 1.033954614161404
$ws.Range("K11").Value = 1.041501651352399
$ws.Range("L11").Value = 1.031135811890018
$ws.Range("M11").Value = 1.051870684095863
$ws.Range("N11").Value = 1.015323472573339

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027141925908321
$ws.Range("D12").Value = 1.037792932542967
$ws.Range("E12").Value = 1.027378734681419
$ws.Range("F12").Value = 1.048158189213696
$ws.Range("I12").Value = 1.039135261391278
$ws.Range("J12").Value = 1.033853139971928
$ws.Range("K12").Value = 1.041398287169909
$ws.Range("L12").Value = 1.031023503928236
$ws.Range("M12").Value = 1.051725174576079
$ws.Range("N12").Value = 1.015289689987763

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027182403565472
$ws.Range("D13").Value = 1.037824845359911
$ws.Range("E13").Value = 1.027412685100272
$ws.Range("F13").Value = 1.04819905983787
$ws.Range("I13").Value = 1.039146012903262
$ws.Range("J13").Value = 1.033874906453436
$ws.Range("K13").Value = 1.041420460062402
$ws.Range("L13").Value = 1.031047591763916
$ws.Range("M13").Value = 1.051756384194477
$ws.Range("N13").Value = 1.01529693660313

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027315054326795
$ws.Range("D14").Value = 1.037929429224594
$ws.Range("E14").Value = 1.027523958507526
$ws.Range("F14").Value = 1.048333007101332
$ws.Range("I14").Value = 1.03918121067045
$ws.Range("J14").Value = 1.033946226240694
$ws.Range("K14").Value = 1.041493107655698
$ws.Range("L14").Value = 1.031126527307449
$ws.Range("M14").Value = 1.051858655033753
$ws.Range("N14").Value = 1.015320680156326

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027396802188113
$ws.Range("D15").Value = 1.037993881417266
$ws.Range("E15").Value = 1.027592542442077
$ws.Range("F15").Value = 1.048415560599519
$ws.Range("I15").Value = 1.039202873922424
$ws.Range("J15").Value = 1.033990168916512
$ws.Range("K15").Value = 1.041537865503632
$ws.Range("L15").Value = 1.031175169698641
$ws.Range("M15").Value = 1.051921675252494
$ws.Range("N15").Value = 1.015335308936338

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027872873505037
$ws.Range("D16").Value = 1.038369242785352
$ws.Range("E16").Value = 1.027992106427434
$ws.Range("F16").Value = 1.048896424959662
$ws.Range("I16").Value = 1.039328607659621
$ws.Range("J16").Value = 1.034245936325626
$ws.Range("K16").Value = 1.041798332711231
$ws.Range("L16").Value = 1.031458404313657
$ws.Range("M16").Value = 1.052288597536424
$ws.Range("N16").Value = 1.015420448183578

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028171730895274
$ws.Range("D17").Value = 1.03860489072881
$ws.Range("E17").Value = 1.028243071719569
$ws.Range("F17").Value = 1.049198380146764
$ws.Range("I17").Value = 1.039407163039834
$ws.Range("J17").Value = 1.034406373156153
$ws.Range("K17").Value = 1.041961678017052
$ws.Range("L17").Value = 1.031636169857384
$ws.Range("M17").Value = 1.052518860047138
$ws.Range("N17").Value = 1.015473847599361

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028346131771393
$ws.Range("D18").Value = 1.038742409339359
$ws.Range("E18").Value = 1.028389573566965
$ws.Range("F18").Value = 1.049374620699578
$ws.Range("I18").Value = 1.039452869584856
$ws.Range("J18").Value = 1.034499953056808
$ws.Range("K18").Value = 1.042056940086017
$ws.Range("L18").Value = 1.031739893029465
$ws.Range("M18").Value = 1.052653204380466
$ws.Range("N18").Value = 1.015504992190751

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028405611883164
$ws.Range("D19").Value = 1.038789311300015
$ws.Range("E19").Value = 1.028439546870015
$ws.Range("F19").Value = 1.049434733644826
$ws.Range("I19").Value = 1.039468435059842
$ws.Range("J19").Value = 1.034531861311257
$ws.Range("K19").Value = 1.042089419488625
$ws.Range("L19").Value = 1.031775265932773
$ws.Range("M19").Value = 1.05269901839378
$ws.Range("N19").Value = 1.015515611272408

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028139657786847
$ws.Range("D20").Value = 1.038579600767437
$ws.Range("E20").Value = 1.028216133294798
$ws.Range("F20").Value = 1.049165971280785
$ws.Range("I20").Value = 1.039398746522107
$ws.Range("J20").Value = 1.034389159812249
$ws.Range("K20").Value = 1.041944154099303
$ws.Range("L20").Value = 1.031617093610194
$ws.Range("M20").Value = 1.052494151310773
$ws.Range("N20").Value = 1.015468118590555

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027275988366052
$ws.Range("D21").Value = 1.03789862883343
$ws.Range("E21").Value = 1.027491186109934
$ws.Range("F21").Value = 1.048293557931887
$ws.Range("I21").Value = 1.039170850648027
$ws.Range("J21").Value = 1.033925224296618
$ws.Range("K21").Value = 1.041471715316683
$ws.Range("L21").Value = 1.03110328117183
$ws.Range("M21").Value = 1.051828537196649
$ws.Range("N21").Value = 1.015313688355568

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026733721048811
$ws.Range("D22").Value = 1.03747111155207
$ws.Range("E22").Value = 1.027036460535244
$ws.Range("F22").Value = 1.047746090446902
$ws.Range("I22").Value = 1.039026547400745
$ws.Range("J22").Value = 1.033633537298318
$ws.Range("K22").Value = 1.041174553837747
$ws.Range("L22").Value = 1.030780557409425
$ws.Range("M22").Value = 1.051410377332244
$ws.Range("N22").Value = 1.015216573785673

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027021114908679
$ws.Range("D23").Value = 1.037697685521284
$ws.Range("E23").Value = 1.027277416349744
$ws.Range("F23").Value = 1.048036212712924
$ws.Range("I23").Value = 1.03910314128974
$ws.Range("J23").Value = 1.033788164913452
$ws.Range("K23").Value = 1.041332095693078
$ws.Range("L23").Value = 1.030951607551904
$ws.Range("M23").Value = 1.051632019112197
$ws.Range("N23").Value = 1.01526805761592

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028154149989576
$ws.Range("D24").Value = 1.03859102799899
$ws.Range("E24").Value = 1.028228305244961
$ws.Range("F24").Value = 1.049180615096649
$ws.Range("I24").Value = 1.039402549934932
$ws.Range("J24").Value = 1.034396937782181
$ws.Range("K24").Value = 1.041952072448372
$ws.Range("L24").Value = 1.031625713235567
$ws.Range("M24").Value = 1.052505316013846
$ws.Range("N24").Value = 1.015470707290463

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029472200125375
$ws.Range("D25").Value = 1.039630407637652
$ws.Range("E25").Value = 1.029336381875763
$ws.Range("F25").Value = 1.050513141275645
$ws.Range("I25").Value = 1.039745542149039
$ws.Range("J25").Value = 1.035103372497636
$ws.Range("K25").Value = 1.042670947497102
$ws.Range("L25").Value = 1.032409363994744
$ws.Range("M25").Value = 1.053520138081361
$ws.Range("N25").Value = 1.015705775669382
